$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (rows 3-6) gets a copy of column E's values from rows 2-5
# (same text, shifted down one row), matching the wrapped-text style
# already used by column E.
$ws.Range("E2:E5").Copy()
$ws.Range("F3").PasteSpecial(-4104)  # xlPasteAll
$excel.CutCopyMode = $false

# Make sure the pasted cells keep the wrapped-text look used elsewhere
# in the sheet (style index 1 in the original file).
$ws.Range("F3:F6").WrapText = $true

# Scroll the view up a bit and select the newly-filled range, matching
# the author's final on-screen selection.
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("F3:F6").Select()
